# Swap the presentation's theme color palette from the "Integral" (Red
# Violet) scheme over to the stock "Office Theme" (Office) scheme — this is
# what the commit did to ppt/theme/theme1.xml (the theme used by the real
# slide master / the slides people actually see).
#
# theme1.xml's <a:clrScheme> carries 12 colors in a fixed order:
#   dk1, lt1, dk2, lt2, accent1..accent6, hlink, folHlink
# That's exactly the order PowerPoint exposes through
# Slide.ThemeColorScheme.Item(1..12), so we drive it through there.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function HexToCOMRGB($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette = Office Theme ("Office" color scheme), in
# dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink order.
$officeTheme = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

for ($i = 1; $i -le $officeTheme.Count; $i++) {
    $tcs.Item($i).RGB = HexToCOMRGB($officeTheme[$i - 1])
}
